$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect Price column (D) updates from Excel's automatic numeric coercion
# by temporarily marking the target cells as Text before assigning the
# literal display strings, then clearing the temporary format again so
# the cells end up with no explicit style (matching the original file).
$priceCells = @(
    "D2",
    "D3",
    "D5",
    "D6",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D16",
    "D17",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D27",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = "44.936.00"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "2.264.96"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("E4").Value = "  -0.78%  "
$ws.Range("D5").Value = "300.39"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").Value = "94.16"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").Value = "0.508"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").Value = "34.03"
$ws.Range("E10").Value = "  -2.19%  "
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "7.20"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("D14").Value = "2.609.11"
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").Value = "2.274.78"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "13.59"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").Value = "0.800"
$ws.Range("E17").Value = "  -3.79%  "
$ws.Range("D18").Value = "44.830.14"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "13.37"
$ws.Range("E19").Value = "  +12.52%  "
$ws.Range("D20").Value = "0.0₃0915"
$ws.Range("E20").Value = "  -2.93%  "
$ws.Range("D21").Value = "6.02"
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("D22").Value = "65.47"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").Value = "238.99"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "2.86"
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").Value = "  -3.34%  "
$ws.Range("D27").Value = "41.25"
$ws.Range("E27").Value = "  +10.87%  "
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("D29").Value = "9.55"
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").Value = "19.58"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").Value = "152.54"
$ws.Range("E31").Value = "  +1.57%  "
$ws.Range("D32").Value = "5.50"
$ws.Range("E32").Value = "  -6.49%  "
$ws.Range("D33").Value = "0.0788"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("E34").Value = "  -2.79%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "0.117"
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "2.90"
$ws.Range("E36").Value = "  -4.64%  "
$ws.Range("D37").Value = "0.102"
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -5.03%  "
$ws.Range("D39").Value = "3.87"
$ws.Range("E39").Value = "  +3.09%  "
$ws.Range("D40").Value = "0.0307"
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("D41").Value = "3.20"
$ws.Range("E41").Value = "  -4.27%  "
$ws.Range("D42").Value = "13.57"
$ws.Range("E42").Value = "  -10.54%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "1.92"
$ws.Range("E44").Value = "  +11.29%  "
$ws.Range("D45").Value = "1.774.47"
$ws.Range("E45").Value = "  -2.75%  "
$ws.Range("D46").Value = "0.192"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("D47").Value = "76.51"
$ws.Range("E47").Value = "  -3.80%  "
$ws.Range("D48").Value = "69.45"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").Value = "95.45"
$ws.Range("E49").Value = "  -2.81%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "7.88"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "53.41"
$ws.Range("E51").Value = "  -0.60%  "

foreach ($cell in $priceCells) {
    $ws.Range($cell).ClearFormats()
}
